$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 70000
$ws.Range("J57").Value = 70000
$ws.Range("L57").Value = 210000
$ws.Range("N57").Value = -210998

$ws.Range("H106").Value = 3633.4614
$ws.Range("I106").Value = 3633.4614
$ws.Range("K106").Value = 3633.4614
$ws.Range("M106").Value = -3002.4614

$ws.Range("H132").Value = 1541.6
$ws.Range("I132").Value = 1152.775
$ws.Range("K132").Value = 3458.325
$ws.Range("M132").Value = -928.3250000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5027.1333
$ws.Range("I61").Value = 2465.889
$ws.Range("J61").Value = 15272.111
$ws.Range("K61").Value = 2465.889
$ws.Range("L61").Value = 15272.111
$ws.Range("M61").Value = -2253.889
$ws.Range("N61").Value = -15696.111

$ws.Range("H102").Value = 1272.1538
$ws.Range("I102").Value = 1268.8
$ws.Range("J102").Value = 1283.3334
$ws.Range("K102").Value = 1268.8
$ws.Range("L102").Value = 1283.3334
$ws.Range("M102").Value = 353.2
$ws.Range("N102").Value = -4527.3334

$ws.Range("H132").Value = 6943.8945
$ws.Range("I132").Value = 5754.92
$ws.Range("K132").Value = 17264.76
$ws.Range("M132").Value = -14734.76

$ws.Range("H136").Value = 5027.1333
$ws.Range("I136").Value = 2465.889
$ws.Range("J136").Value = 15272.111
$ws.Range("K136").Value = 7397.667
$ws.Range("L136").Value = 45816.333
$ws.Range("M136").Value = -4847.667
$ws.Range("N136").Value = -50916.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 27197738
$ws.Range("I86").Value = 8965803
$ws.Range("J86").Value = 55558524
$ws.Range("K86").Value = 8965803
$ws.Range("L86").Value = 55558524
$ws.Range("M86").Value = -8964680
$ws.Range("N86").Value = -55560770

$ws.Range("H89").Value = 27197738
$ws.Range("I89").Value = 8965803
$ws.Range("J89").Value = 55558524
$ws.Range("K89").Value = 44829015
$ws.Range("L89").Value = 277792620
$ws.Range("M89").Value = -44823399
$ws.Range("N89").Value = -277803852

$ws.Range("H107").Value = 34094984
$ws.Range("I107").Value = 41670396
$ws.Range("K107").Value = 41670396
$ws.Range("M107").Value = -41668476

$ws.Range("H134").Value = 5561.6445
$ws.Range("I134").Value = 2268.2173
$ws.Range("J134").Value = 9004.772000000001
$ws.Range("K134").Value = 6804.651899999999
$ws.Range("L134").Value = 27014.316
$ws.Range("M134").Value = -4269.651899999999
$ws.Range("N134").Value = -32084.316

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 141.11111
$ws.Range("I7").Value = 38
$ws.Range("J7").Value = 192.66667
$ws.Range("K7").Value = 38
$ws.Range("L7").Value = 192.66667
$ws.Range("M7").Value = 75
$ws.Range("N7").Value = -418.66667

$ws.Range("H58").Value = 17249930
$ws.Range("I58").Value = 50003544
$ws.Range("J58").Value = 11187.526
$ws.Range("K58").Value = 50003544
$ws.Range("L58").Value = 11187.526
$ws.Range("M58").Value = -50003341
$ws.Range("N58").Value = -11593.526

$ws.Range("H99").Value = 8626.076999999999
$ws.Range("I99").Value = 10481.5
$ws.Range("J99").Value = 7035.7144
$ws.Range("K99").Value = 10481.5
$ws.Range("L99").Value = 7035.7144
$ws.Range("M99").Value = -8983.5
$ws.Range("N99").Value = -10031.7144

$ws.Range("H126").Value = 8626.076999999999
$ws.Range("I126").Value = 10481.5
$ws.Range("J126").Value = 7035.7144
$ws.Range("K126").Value = 31444.5
$ws.Range("L126").Value = 21107.1432
$ws.Range("M126").Value = -28974.5
$ws.Range("N126").Value = -26047.1432

$ws.Range("H132").Value = 6551.857
$ws.Range("I132").Value = 2972.8333
$ws.Range("J132").Value = 9236.125
$ws.Range("K132").Value = 8918.499899999999
$ws.Range("L132").Value = 27708.375
$ws.Range("M132").Value = -6388.499899999999
$ws.Range("N132").Value = -32768.375

$ws.Range("H134").Value = 7807.793
$ws.Range("I134").Value = 3602.2856
$ws.Range("K134").Value = 10806.8568
$ws.Range("M134").Value = -8271.856800000001

$ws.Range("H136").Value = 17249930
$ws.Range("I136").Value = 50003544
$ws.Range("J136").Value = 11187.526
$ws.Range("K136").Value = 150010632
$ws.Range("L136").Value = 33562.578
$ws.Range("M136").Value = -150008082
$ws.Range("N136").Value = -38662.578

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 2200
$ws.Range("I88").Value = 2200
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 6600
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("M88").Value = -6172

$ws.Range("H91").Value = 2200
$ws.Range("I91").Value = 2200
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 6600
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("M91").Value = -5118

$ws.Range("H138").Value = 53975.7
$ws.Range("I138").Value = 75066.36
$ws.Range("K138").Value = 225199.08
$ws.Range("M138").Value = -220059.08

$ws.Range("H139").Value = 32242.945
$ws.Range("I139").Value = 60970.06
$ws.Range("K139").Value = 182910.18
$ws.Range("M139").Value = -177770.18

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5799.273
$ws.Range("I97").Value = 5048.6665
$ws.Range("J97").Value = 6700
$ws.Range("K97").Value = 5048.6665
$ws.Range("L97").Value = 6700
$ws.Range("M97").Value = -4552.6665
$ws.Range("N97").Value = -7692

$ws.Range("H102").Value = 8732.538
$ws.Range("I102").Value = 7752.8
$ws.Range("K102").Value = 7752.8
$ws.Range("M102").Value = -6130.8

$ws.Range("H122").Value = 1543794.9
$ws.Range("I122").Value = 2071662.4
$ws.Range("J122").Value = 4181.4165
$ws.Range("K122").Value = 6214987.199999999
$ws.Range("L122").Value = 12544.2495
$ws.Range("M122").Value = -6212537.199999999
$ws.Range("N122").Value = -17444.2495

$ws.Range("H132").Value = 8719.583000000001
$ws.Range("I132").Value = 2865.5
$ws.Range("K132").Value = 8596.5
$ws.Range("M132").Value = -6066.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7432.706
$ws.Range("I7").Value = 6634.7
$ws.Range("K7").Value = 6634.7
$ws.Range("M7").Value = -6522.7

$ws.Range("H46").Value = 3549.6365
$ws.Range("J46").Value = 3805.75
$ws.Range("L46").Value = 3805.75
$ws.Range("N46").Value = -4181.75

$ws.Range("H100").Value = 3216.4546
$ws.Range("I100").Value = 3102.5715
$ws.Range("J100").Value = 3269.6
$ws.Range("K100").Value = 3102.5715
$ws.Range("L100").Value = 3269.6
$ws.Range("M100").Value = -2561.5715
$ws.Range("N100").Value = -4351.6

$ws.Range("H126").Value = 7432.706
$ws.Range("I126").Value = 6634.7
$ws.Range("K126").Value = 19904.1
$ws.Range("M126").Value = -17434.1

$ws.Range("H132").Value = 16675063
$ws.Range("I132").Value = 41669824
$ws.Range("J132").Value = 11888.889
$ws.Range("K132").Value = 125009472
$ws.Range("L132").Value = 35666.667
$ws.Range("M132").Value = -125006942
$ws.Range("N132").Value = -40726.667

$ws.Range("H136").Value = 10529.042
$ws.Range("I136").Value = 2385.2856
$ws.Range("J136").Value = 13882.353
$ws.Range("K136").Value = 7155.8568
$ws.Range("L136").Value = 41647.05899999999
$ws.Range("M136").Value = -4605.8568
$ws.Range("N136").Value = -46747.05899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1710.8667
$ws.Range("I96").Value = 1815.6666
$ws.Range("K96").Value = 1815.6666
$ws.Range("M96").Value = -442.6666

$ws.Range("H100").Value = 936.4
$ws.Range("I100").Value = 589.8
$ws.Range("J100").Value = 1629.6
$ws.Range("K100").Value = 1179.6
$ws.Range("L100").Value = 3259.2
$ws.Range("M100").Value = -638.5999999999999
$ws.Range("N100").Value = -4341.2

$ws.Range("H132").Value = 13517947
$ws.Range("I132").Value = 17861816
$ws.Range("K132").Value = 53585448
$ws.Range("M132").Value = -53582918

